$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "27.526.71"
$ws.Cells.Item(2, 5).Value = "  -3.05%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.757.04"
$ws.Cells.Item(3, 5).Value = "  -2.69%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.007"
$ws.Cells.Item(4, 5).Value = "  +0.58%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "324.93"
$ws.Cells.Item(5, 5).Value = "  -0.79%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.003"
$ws.Cells.Item(6, 5).Value = "  +0.43%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4468"
$ws.Cells.Item(7, 5).Value = "  +0.33%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3695"
$ws.Cells.Item(8, 5).Value = "  -1.10%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "45.13"
$ws.Cells.Item(9, 5).Value = "  +0.83%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.07695"
$ws.Cells.Item(10, 5).Value = "  +2.40%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "1.116"
$ws.Cells.Item(11, 5).Value = "  -2.80%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.005"
$ws.Cells.Item(12, 5).Value = "  +0.44%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "21.61"
$ws.Cells.Item(13, 5).Value = "  -4.22%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.153"
$ws.Cells.Item(14, 5).Value = "  -2.30%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "7.372"
$ws.Cells.Item(15, 5).Value = "  -4.29%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "1.762.61"
$ws.Cells.Item(16, 5).Value = "  -2.02%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "90.46"
$ws.Cells.Item(17, 5).Value = "  +11.88%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.00001071"
$ws.Cells.Item(18, 5).Value = "  -2.10%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06278"
$ws.Cells.Item(19, 5).Value = "  -7.53%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "1.003"
$ws.Cells.Item(20, 5).Value = "  +0.40%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "17.38"
$ws.Cells.Item(21, 5).Value = "  -0.51%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.155"
$ws.Cells.Item(22, 5).Value = "  -2.71%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.5326"
$ws.Cells.Item(23, 5).Value = "  -2.37%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "27.563.94"
$ws.Cells.Item(24, 5).Value = "  -2.81%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "11.55"
$ws.Cells.Item(25, 5).Value = "  -2.19%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.319"
$ws.Cells.Item(26, 5).Value = "  -3.79%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "20.56"
$ws.Cells.Item(27, 5).Value = "  +0.39%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "152.98"
$ws.Cells.Item(28, 5).Value = "  -0.30%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.297"
$ws.Cells.Item(29, 5).Value = "  -2.36%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.961.24"
$ws.Cells.Item(30, 5).Value = "  -2.09%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "127.55"
$ws.Cells.Item(31, 5).Value = "  -3.83%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.182"
$ws.Cells.Item(32, 5).Value = "  -5.76%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "5.712"
$ws.Cells.Item(33, 5).Value = "  -1.92%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.09211"
$ws.Cells.Item(34, 5).Value = "  -1.34%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.648"
$ws.Cells.Item(35, 5).Value = "  -9.01%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "12.59"
$ws.Cells.Item(36, 5).Value = "  +3.88%  "
$ws.Cells.Item(37, 2).Value = "Algorand"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.2164"
$ws.Cells.Item(37, 5).Value = "  -5.42%  "
$ws.Cells.Item(38, 2).Value = "VeChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.02310"
$ws.Cells.Item(38, 5).Value = "  -0.47%  "
$ws.Cells.Item(39, 2).Value = "Hedera"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.06080"
$ws.Cells.Item(39, 5).Value = "  -4.13%  "
$ws.Cells.Item(40, 2).Value = "TheSandbox"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.6431"
$ws.Cells.Item(40, 5).Value = "  -2.17%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "5.040"
$ws.Cells.Item(41, 5).Value = "  -2.42%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.173"
$ws.Cells.Item(42, 5).Value = "  -3.22%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "7.952"
$ws.Cells.Item(43, 5).Value = "  -2.78%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.004"
$ws.Cells.Item(44, 5).Value = "  +0.46%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.403"
$ws.Cells.Item(45, 5).Value = "  -3.93%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "13.62"
$ws.Cells.Item(46, 5).Value = "  -3.24%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.5961"
$ws.Cells.Item(47, 5).Value = "  -1.77%  "
$ws.Cells.Item(48, 5).Value = "  -1.51%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "125.52"
$ws.Cells.Item(49, 5).Value = "  -2.18%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.993"
$ws.Cells.Item(50, 5).Value = "  -2.00%  "
$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.06897"
$ws.Cells.Item(51, 5).Value = "  -2.97%  "
